$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = -3
    4  = -5
    5  = -3
    6  = -1
    7  = -5
    8  = -1
    9  = -2
    10 = -3
    12 = -1
    13 = 1
    14 = -1
    16 = -2
    17 = -7
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
